# Update the products export sheet with the 31st-edition brand data.
# Row 2 (existing) is overwritten in place and two new rows (3 and 4) are
# appended, all for the same brand ("B10" / code "0277" / WHISKY) but with
# different barcodes, pack sizes (brand_size) and prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# brand_code / barcode / brand_size look numeric (leading zeros, long digit
# strings) but must stay text, exactly like the original export - so mark
# those columns as Text before writing into them, which keeps Excel from
# silently converting "0277" -> 277 or mangling the 13-digit barcodes.
$ws.Range("B2:C4").NumberFormat = "@"
$ws.Range("H2:H4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "B10"
$ws.Range("B2").Value = "0277"
$ws.Range("C2").Value = "8901544027077"
$ws.Range("D2").Value = 320
$ws.Range("E2").Value = 320
$ws.Range("F2").Value = 370
$ws.Range("G2").Value = "WHISKY"
$ws.Range("H2").Value = "180"
$ws.Range("I2").Value = "static/uploads/Default.png"

# Row 3 (new)
$ws.Range("A3").Value = "B10"
$ws.Range("B3").Value = "0277"
$ws.Range("C3").Value = "8901544027060"
$ws.Range("D3").Value = 730
$ws.Range("E3").Value = 730
$ws.Range("F3").Value = 740
$ws.Range("G3").Value = "WHISKY"
$ws.Range("H3").Value = "375"
$ws.Range("I3").Value = "static/uploads/Default.png"

# Row 4 (new)
$ws.Range("A4").Value = "B10"
$ws.Range("B4").Value = "0277"
$ws.Range("C4").Value = "8901544027039"
$ws.Range("D4").Value = 1450
$ws.Range("E4").Value = 1450
$ws.Range("F4").Value = 1480
$ws.Range("G4").Value = "WHISKY"
$ws.Range("H4").Value = "750"
$ws.Range("I4").Value = "static/uploads/Default.png"

Write-Host "Updated rows 2-4 of Sheet1 with new brand data"
